$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Component Label for existing 5V Regulator part ---
$ws.Range("I3").Value = "U1"

# --- Row 4: Nichicon 0.33uF Electrolytic Capacitor ---
$ws.Range("A4").Value = "Nichicon 0.33uF Electrolytic Capacitor"
$ws.Range("B4").Value = "0.33uF 50V 4mm Radial Through Hole"
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0.076
$ws.Range("F4").Value = "RS"
$ws.Range("G4").Value = "475-8983"
$ws.Range("H4").Value = "USR1HR33MDD"
$ws.Range("I4").Value = "C1"

# --- Row 5: Nichicon 1uF Electrolytic Capacitor ---
$ws.Range("A5").Value = "Nichicon 1uF Electrolytic Capacitor"
$ws.Range("B5").Value = "1uF 50V 5mm Radial Through Hole"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 0.111
$ws.Range("F5").Value = "RS"
$ws.Range("G5").Value = "715-2808"
$ws.Range("H5").Value = "UPW1H010MDD"
$ws.Range("I5").Value = "C2"

# --- Row 6: Mercury 4MHz Crystal Oscillator ---
$ws.Range("A6").Value = "Mercury 4MHz Crystal Oscillator"
$ws.Range("B6").Value = "4MHz +-50ppm HCMOS TTL Crystal Oscillator"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1.81
$ws.Range("F6").Value = "RS"
$ws.Range("G6").Value = "767-5244"
$ws.Range("H6").Value = "5H14ET-4.000"
$ws.Range("I6").Value = "O1"

# Move the active selection as recorded in the edited workbook
$ws.Range("H24").Select()
